# Piaps prak 4 added
# Re-splits three list-item paragraphs into several <w:r> runs (mirroring
# the incremental hand-edits the author made directly in Word) and adds a
# yellow highlight on part of the "requirements" bullet.

$d = $word.ActiveDocument

function Replace-Paragraph($findText, $pTagAndPPr, $runsXml) {
    $target = $d.Content.Duplicate
    $target.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $ins = $target.Duplicate
    $ins.Collapse(0)

    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body>' + $pTagAndPPr + $runsXml + '</w:p></w:body></w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'

    $ins.InsertXML($xml)
}

$rPrPlain = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr>'
$rPrHighlight = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:highlight w:val="yellow"/></w:rPr>'

# ---------------------------------------------------------------------
# Bullet 1: "Поддержка сценариев аварий и сбоев." ->
#           "Поддержка сценариев при авариях и сбоях."
# ---------------------------------------------------------------------
$p1Tag = '<w:p w14:paraId="001F0E3C" w14:textId="77777777" w:rsidR="007B7427" w:rsidRPr="007B7427" w:rsidRDefault="007B7427" w:rsidP="007B7427">' +
         '<w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:before="74"/><w:ind w:right="815"/>' +
         '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr></w:pPr>'

$p1Runs  = '<w:r w:rsidRPr="007B7427">' + $rPrPlain + '<w:t>Поддержка сценариев</w:t></w:r>'
$p1Runs += '<w:r>' + $rPrPlain + '<w:t xml:space="preserve"> </w:t></w:r>'
$p1Runs += '<w:r>' + $rPrPlain + '<w:t>при</w:t></w:r>'
$p1Runs += '<w:r>' + $rPrPlain + '<w:t xml:space="preserve"> авари</w:t></w:r>'
$p1Runs += '<w:r>' + $rPrPlain + '<w:t>ях</w:t></w:r>'
$p1Runs += '<w:r>' + $rPrPlain + '<w:t xml:space="preserve"> и сбо</w:t></w:r>'
$p1Runs += '<w:r>' + $rPrPlain + '<w:t>ях</w:t></w:r>'
$p1Runs += '<w:r>' + $rPrPlain + '<w:t>.</w:t></w:r>'

Replace-Paragraph "Поддержка сценариев аварий и сбоев." $p1Tag $p1Runs

# ---------------------------------------------------------------------
# Bullet 2: "Составление перечня ключевых функций и требований
#            безопасности к системе." -> same text, "безопасности к
#            системе" highlighted yellow.
# ---------------------------------------------------------------------
$p2Tag = '<w:p w14:paraId="346E902D" w14:textId="77777777" w:rsidR="00C905DA" w:rsidRPr="00A265B6" w:rsidRDefault="00C905DA" w:rsidP="00A265B6">' +
         '<w:pPr><w:pStyle w:val="a7"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:spacing w:before="74"/><w:ind w:right="815"/>' +
         '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr></w:pPr>'

$p2Runs  = '<w:r w:rsidRPr="00A265B6">' + $rPrPlain + '<w:t xml:space="preserve">Составление перечня ключевых функций и требований </w:t></w:r>'
$p2Runs += '<w:r>' + $rPrHighlight + '<w:t>безопасности к системе</w:t></w:r>'
$p2Runs += '<w:r>' + $rPrPlain + '<w:t>.</w:t></w:r>'

Replace-Paragraph "Составление перечня ключевых функций и требований безопасности к системе." $p2Tag $p2Runs

# ---------------------------------------------------------------------
# Bullet 3: "Разработка моделей технических процессов и сценариев для
#            моделирования рисков." -> "... рисков аварий (сбоев)."
# ---------------------------------------------------------------------
$p3Tag = '<w:p w14:paraId="242811E5" w14:textId="77777777" w:rsidR="00C905DA" w:rsidRPr="00AB6F16" w:rsidRDefault="00C905DA" w:rsidP="00AB6F16">' +
         '<w:pPr><w:pStyle w:val="a7"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:spacing w:before="74"/><w:ind w:right="815"/>' +
         '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr></w:pPr>'

$p3Runs  = '<w:r w:rsidRPr="00AB6F16">' + $rPrPlain + '<w:t>Разработка моделей технических процессов и сценариев для моделирования рисков</w:t></w:r>'
$p3Runs += '<w:r>' + $rPrPlain + '<w:t xml:space="preserve"> аварий (сбоев)</w:t></w:r>'
$p3Runs += '<w:r>' + $rPrPlain + '<w:t>.</w:t></w:r>'

Replace-Paragraph "Разработка моделей технических процессов и сценариев для моделирования рисков." $p3Tag $p3Runs

Write-Output "done"
